$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows per repull/push of data
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -3
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -9
$ws.Range("F12").Value = 0
$ws.Range("F17").Value = -3
